$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new "American Political Science Review" research-note row
# (row 91) to the bottom of the table on Sheet1.
# Shared-string indices are assigned in first-seen order, so write the
# "7k words" (index 208) value before the link text (index 209) to keep
# the dedup table in the same order as the canonical workbook.
$ws.Range("A91").Value = "American Political Science Review"
$ws.Range("B91").Value = 82
$ws.Range("D91").Value = "7k words"
$ws.Range("C91").Value = "<a href='https://www.cambridge.org/core/journals/american-political-science-review/information/author-instructions/preparing-your-materials'target='_blank'>Research Note</a>"

# Match the author's final selection/view state.
$ws.Range("A2").Select()
